$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 524.3889
$ws.Range("J6").Value = 1020.4
$ws.Range("L6").Value = 3061.2
$ws.Range("N6").Value = -3285.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2750
$ws.Range("I43").Value = 2000
$ws.Range("K43").Value = 2000
$ws.Range("M43").Value = -1931

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 443.07144
$ws.Range("I55").Value = 423.66666
$ws.Range("J55").Value = 478
$ws.Range("K55").Value = 423.66666
$ws.Range("L55").Value = 478
$ws.Range("M55").Value = -209.66666
$ws.Range("N55").Value = -906

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4653.091
$ws.Range("I113").Value = 4595.625
$ws.Range("K113").Value = 4595.625
$ws.Range("M113").Value = -1341.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2201.1
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("M116").Value = 1442

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12951.218
$ws.Range("I132").Value = 11440.947
$ws.Range("K132").Value = 34322.841
$ws.Range("M132").Value = -31792.841

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1552.8889
$ws.Range("I137").Value = 996.5714
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 2989.7142
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -439.7142000000003
$ws.Range("N137").Value = -15600

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3226.7917
$ws.Range("I138").Value = 883.3
$ws.Range("K138").Value = 2649.9
$ws.Range("M138").Value = 2490.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2591.5
$ws.Range("I3").Value = 600
$ws.Range("J3").Value = 4583
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 4583
$ws.Range("M3").Value = -485
$ws.Range("N3").Value = -4813

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 85.75
$ws.Range("I4").Value = 90
$ws.Range("K4").Value = 90
$ws.Range("M4").Value = 26

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 159.8
$ws.Range("I5").Value = 199.66667
$ws.Range("K5").Value = 199.66667
$ws.Range("M5").Value = -87.66667000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1953.3334
$ws.Range("I45").Value = 1123.2222
$ws.Range("K45").Value = 1123.2222
$ws.Range("M45").Value = -746.2221999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 11900.5
$ws.Range("I63").Value = 2235
$ws.Range("J63").Value = 17699.8
$ws.Range("K63").Value = 2235
$ws.Range("L63").Value = 17699.8
$ws.Range("M63").Value = -1549
$ws.Range("N63").Value = -19071.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 11900.5
$ws.Range("I66").Value = 2235
$ws.Range("J66").Value = 17699.8
$ws.Range("K66").Value = 11175
$ws.Range("L66").Value = 88499
$ws.Range("M66").Value = -7743
$ws.Range("N66").Value = -95363

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 796.1539
$ws.Range("I97").Value = 444.91666
$ws.Range("K97").Value = 444.91666
$ws.Range("M97").Value = 51.08334000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 159.8
$ws.Range("I4").Value = 199.66667
$ws.Range("K4").Value = 199.66667
$ws.Range("M4").Value = -84.66667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 30166.666
$ws.Range("I26").Value = 30166.666
$ws.Range("K26").Value = 30166.666
$ws.Range("M26").Value = -29874.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 539.5714
$ws.Range("I94").Value = 296.33334
$ws.Range("J94").Value = 1999
$ws.Range("K94").Value = 296.33334
$ws.Range("L94").Value = 1999
$ws.Range("M94").Value = 154.66666
$ws.Range("N94").Value = -2901

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1669.6666
$ws.Range("I105").Value = 1505
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 1505
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = 242
$ws.Range("N105").Value = -5493

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2448.3438
$ws.Range("I107").Value = 1667.138
$ws.Range("K107").Value = 1667.138
$ws.Range("M107").Value = 252.8620000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20000
$ws.Range("I4").Value = 20000
$ws.Range("K4").Value = 20000
$ws.Range("M4").Value = -19888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3433.6775
$ws.Range("I7").Value = 5993.4116
$ws.Range("J7").Value = 325.42856
$ws.Range("K7").Value = 5993.4116
$ws.Range("L7").Value = 325.42856
$ws.Range("M7").Value = -5880.4116
$ws.Range("N7").Value = -551.4285600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2519.6
$ws.Range("I22").Value = 2024.5
$ws.Range("K22").Value = 2024.5
$ws.Range("M22").Value = -1674.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3961.0488
$ws.Range("I31").Value = 1914.9656
$ws.Range("K31").Value = 1914.9656
$ws.Range("M31").Value = -1619.9656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3961.0488
$ws.Range("I34").Value = 1914.9656
$ws.Range("K34").Value = 1914.9656
$ws.Range("M34").Value = -1712.9656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4268.091
$ws.Range("I58").Value = 3334.8333
$ws.Range("J58").Value = 5388
$ws.Range("K58").Value = 3334.8333
$ws.Range("L58").Value = 5388
$ws.Range("M58").Value = -3131.8333
$ws.Range("N58").Value = -5794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5749
$ws.Range("J62").Value = 4872.5
$ws.Range("L62").Value = 4872.5
$ws.Range("N62").Value = -6120.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5749
$ws.Range("J65").Value = 4872.5
$ws.Range("L65").Value = 24362.5
$ws.Range("N65").Value = -30602.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 53375
$ws.Range("J68").Value = 62833.332
$ws.Range("L68").Value = 62833.332
$ws.Range("N68").Value = -64331.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 53375
$ws.Range("J71").Value = 62833.332
$ws.Range("L71").Value = 188499.996
$ws.Range("N71").Value = -195987.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 65382.5
$ws.Range("J74").Value = 65382.5
$ws.Range("L74").Value = 65382.5
$ws.Range("N74").Value = -67130.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 65382.5
$ws.Range("J77").Value = 65382.5
$ws.Range("L77").Value = 196147.5
$ws.Range("N77").Value = -204883.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -42620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4268.091
$ws.Range("I136").Value = 3334.8333
$ws.Range("J136").Value = 5388
$ws.Range("K136").Value = 10004.4999
$ws.Range("L136").Value = 16164
$ws.Range("M136").Value = -7454.499899999999
$ws.Range("N136").Value = -21264

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1595.2727
$ws.Range("J5").Value = 1674.6666
$ws.Range("L5").Value = 5023.9998
$ws.Range("N5").Value = -5247.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3724.5833
$ws.Range("I132").Value = 2199.5
$ws.Range("J132").Value = 4487.125
$ws.Range("K132").Value = 19795.5
$ws.Range("L132").Value = 40384.125
$ws.Range("M132").Value = -17265.5
$ws.Range("N132").Value = -45444.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1595.2727
$ws.Range("J135").Value = 1674.6666
$ws.Range("L135").Value = 15071.9994
$ws.Range("N135").Value = -20141.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 397.57144
$ws.Range("I2").Value = 273.6
$ws.Range("J2").Value = 466.44446
$ws.Range("K2").Value = 273.6
$ws.Range("L2").Value = 466.44446
$ws.Range("M2").Value = -160.6
$ws.Range("N2").Value = -692.4444599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 17337732
$ws.Range("I3").Value = 23503472
$ws.Range("J3").Value = 5006250
$ws.Range("K3").Value = 23503472
$ws.Range("L3").Value = 5006250
$ws.Range("M3").Value = -23503356
$ws.Range("N3").Value = -5006482

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4599.6
$ws.Range("I70").Value = 4499.75
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 4499.75
$ws.Range("L70").Value = 4999
$ws.Range("M70").Value = -4229.75
$ws.Range("N70").Value = -5539

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4599.6
$ws.Range("I73").Value = 4499.75
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 4499.75
$ws.Range("L73").Value = 4999
$ws.Range("M73").Value = -3563.75
$ws.Range("N73").Value = -6871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3497
$ws.Range("I80").Value = 3495
$ws.Range("J80").Value = 3498.3333
$ws.Range("K80").Value = 3495
$ws.Range("L80").Value = 3498.3333
$ws.Range("M80").Value = -2497
$ws.Range("N80").Value = -5494.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3497
$ws.Range("I83").Value = 3495
$ws.Range("J83").Value = 3498.3333
$ws.Range("K83").Value = 17475
$ws.Range("L83").Value = 17491.6665
$ws.Range("M83").Value = -12483
$ws.Range("N83").Value = -27475.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6301.385
$ws.Range("I7").Value = 6165.222
$ws.Range("K7").Value = 6165.222
$ws.Range("M7").Value = -6053.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6301.385
$ws.Range("I126").Value = 6165.222
$ws.Range("K126").Value = 18495.666
$ws.Range("M126").Value = -16025.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20000000
$ws.Range("I5").Value = 20000000
$ws.Range("J5").Value = 20000000
$ws.Range("K5").Value = 20000000
$ws.Range("L5").Value = 20000000
$ws.Range("M5").Value = -19999888
$ws.Range("N5").Value = -20000224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4844.7896
$ws.Range("I126").Value = 2794.6667
$ws.Range("K126").Value = 8384.000100000001
$ws.Range("M126").Value = -5914.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1813.5555
$ws.Range("I132").Value = 1639.7858
$ws.Range("K132").Value = 4919.357400000001
$ws.Range("M132").Value = -2389.357400000001
